# Generate Report for Handback
# Updates the handback status report with freshly generated timestamps
# (and the handoff type for the 4d9002ec entry) to reflect the latest
# localization generation run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# "Latest HO Xliff Generate Date" for 4d9002ec... and 9ed00e70... on the
# Overview sheet (both rows shared the same generation timestamp).
$wsOverview.Range("G3").Value = "2016-08-28 04:17:27"
$wsOverview.Range("G4").Value = "2016-08-28 04:17:27"

# zh-cn sheet: Priority/handoff-type for the 4d9002ec row changed from
# "ht" to "mt", and its handoff/handback datetimes were refreshed.
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("E4").Value = "mt"
$wsZhCn.Range("H3").Value = "2016-08-28 04:17:23"
$wsZhCn.Range("H4").Value = "2016-08-28 04:17:23"
$wsZhCn.Range("K3").Value = "2016-08-28 04:17:39"
$wsZhCn.Range("K4").Value = "2016-08-28 04:17:39"

# de-de sheet: Correspond Handoff Datetime for the 4d9002ec / 9ed00e70
# rows (shared "Latest HO Xliff Generate Date" with the Overview sheet)
# and the Correspond Handback DateTime were refreshed.
$wsDeDe.Range("H3").Value = "2016-08-28 04:17:27"
$wsDeDe.Range("H4").Value = "2016-08-28 04:17:27"
$wsDeDe.Range("K3").Value = "2016-08-28 04:17:45"
$wsDeDe.Range("K4").Value = "2016-08-28 04:17:45"
